$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (29) down onto the three
# new rows (30-32) so the new cells pick up the same cell styles (s="1" for
# appid, s="2" for the email/recovery columns, etc.) as the rest of the table.
$ws.Range("A29:G29").Copy()
$ws.Range("A30:G32").PasteSpecial(-4122)

# Row 30
$ws.Range("A30").Value = "com.hamxa.shaynachim"
$ws.Range("B30").Value = "bitcoin free"
$ws.Range("C30").Value = "mirogya57@gmail.com"
$ws.Range("D30").Value = "Mirogyagi767@live.com"
$ws.Range("E30").Value = "27/5/2019 15:59"
$ws.Range("F30").Value = "This innovative app explains as easy as it can what is bitcoin. Really good"
$ws.Range("G30").Value = "no"

# Row 31
$ws.Range("A31").Value = "com.hamxa.shaynachim"
$ws.Range("B31").Value = "bitcoin free"
$ws.Range("C31").Value = "galiatia942@gmail.com"
$ws.Range("D31").Value = "syechimovitz@gmail.com"
$ws.Range("E31").Value = "27/5/2019 15:59"
$ws.Range("F31").Value = "free app and free info about bitcoin. Splendid!!!"
$ws.Range("G31").Value = "no"

# Row 32
$ws.Range("A32").Value = "com.hamxa.shaynachim"
$ws.Range("B32").Value = "bitcoin free"
$ws.Range("C32").Value = "irisalmog47@gmail.com"
$ws.Range("D32").Value = "bittonnir12@gmail.com"
$ws.Range("E32").Value = "27/5/2019 15:59"
$ws.Range("F32").Value = "very special app info. Keep it simple and really educative"
$ws.Range("G32").Value = "no"

# New hyperlink on D32, matching the existing mailto: hyperlinks in the sheet.
# Hyperlinks.Add applies Excel's default blue/underlined "Hyperlink" style to
# the target cell; re-paste the column's normal formatting afterwards so D32
# keeps the same style as the rest of the table (matching the other
# hyperlinked cells such as D2, D27, ...).
$ws.Hyperlinks.Add($ws.Range("D32"), "mailto:bittonnir12@gmail.com", "", "", "bittonnir12@gmail.com")
$ws.Range("D31").Copy()
$ws.Range("D32").PasteSpecial(-4122)

# Scroll / selection state, as recorded in the saved workbook.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F26:F32").Select()
